$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.74998539686203
$ws.Range("B1").Value = 1.570627927780151
$ws.Range("C1").Value = 4.555920600891113
$ws.Range("D1").Value = 2.400482416152954
$ws.Range("E1").Value = 1.261965751647949
